$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("臨床イベント")

# Column A holds dates stored as text (matches source file's
# numberStoredAsText convention) - force text formatting on the rows
# whose date value changes, so the new date-like strings aren't
# auto-converted to date serials. Rows 2-3 keep their original
# (already-text) date value untouched.
$ws.Range("A4:A8").NumberFormat = "@"

# Row 2 (date unchanged - 2025-01-01)
$ws.Cells.Item(2, 2).Value = "低血糖"
$ws.Cells.Item(2, 3).Value = "血糖28mg/dL、哺乳不良あり"

# Row 3 (date unchanged - 2025-01-01)
$ws.Cells.Item(3, 2).Value = "意識障害"
$ws.Cells.Item(3, 3).Value = "傾眠傾向、刺激で開眼"

# Row 4
$ws.Cells.Item(4, 1).Value = "2025-01-01"
$ws.Cells.Item(4, 2).Value = "てんかん発作"
$ws.Cells.Item(4, 3).Value = "強直間代発作1回、約2分間"

# Row 5
$ws.Cells.Item(5, 1).Value = "2025-01-02"
$ws.Cells.Item(5, 2).Value = "低血糖"
$ws.Cells.Item(5, 3).Value = "血糖32mg/dL、GIR増量後"

# Row 6
$ws.Cells.Item(6, 1).Value = "2025-01-03"
$ws.Cells.Item(6, 2).Value = "低血糖"
$ws.Cells.Item(6, 3).Value = "血糖35mg/dL、ジアゾキシド開始"

# Row 7
$ws.Cells.Item(7, 1).Value = "2025-01-05"
$ws.Cells.Item(7, 2).Value = "低血糖"
$ws.Cells.Item(7, 3).Value = "血糖48mg/dL、改善傾向"

# Row 8 (new)
$ws.Cells.Item(8, 1).Value = "2025-01-07"
$ws.Cells.Item(8, 2).Value = "低血糖"
$ws.Cells.Item(8, 3).Value = "血糖55mg/dL、GIR減量可能に"
